$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Ano" in D1, matching the style/formatting of the other header cells (A1:C1)
$ws.Range("D1").Value = "Ano"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Fill D2:D10 with the reference year 2023
$ws.Range("D2:D10").Value = 2023
